$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("F3").Value = 0.125
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 0.375
$ws.Range("I3").Value = 9
$ws.Range("H5").Value = 0.085714285714285715
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.028571428571428571
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 0.11428571428571428
$ws.Range("O5").Value = 4
$ws.Range("J6").Value = 0.046511627906976744
$ws.Range("K6").Value = 2
$ws.Range("F7").Value = 0.2
$ws.Range("G7").Value = 2
$ws.Range("N7").Value = 0.2
$ws.Range("O7").Value = 2
$ws.Range("L8").Value = 0.096774193548387094
$ws.Range("M8").Value = 3
$ws.Range("F9").Value = 0.064516129032258063
$ws.Range("G9").Value = 2
$ws.Range("L9").Value = 0.032258064516129031
$ws.Range("M9").Value = 1
$ws.Range("H10").Value = 0.032258064516129031
$ws.Range("I10").Value = 1
$ws.Range("N10").Value = 0.096774193548387094
$ws.Range("O10").Value = 3
$ws.Range("F11").Value = 0.086956521739130432
$ws.Range("G11").Value = 2
$ws.Range("J12").Value = 0.064516129032258063
$ws.Range("K12").Value = 2
$ws.Range("F13").Value = 0.17857142857142858
$ws.Range("G13").Value = 5
$ws.Range("L13").Value = 0.071428571428571425
$ws.Range("M13").Value = 2
$ws.Range("C15").Value = 20
$ws.Range("D15").Value = 0.1
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.25
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 0.35
$ws.Range("I15").Value = 7
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("D16").Value = 0.037037037037037035
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.1111111111111111
$ws.Range("G16").Value = 3
$ws.Range("F20").Value = 0.02564102564102564
$ws.Range("G20").Value = 1
$ws.Range("L20").Value = 0.05128205128205128
$ws.Range("M20").Value = 2
$ws.Range("L22").Value = 0.13636363636363635
$ws.Range("M22").Value = 3
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("J29").Value = 0.083333333333333329
$ws.Range("K29").Value = 3
$ws.Range("L29").Value = 0.083333333333333329
$ws.Range("M29").Value = 3
$ws.Range("D32").Value = 0.060606060606060608
$ws.Range("E32").Value = 2
$ws.Range("F32").Value = 0.12121212121212122
$ws.Range("G32").Value = 4
$ws.Range("L32").Value = 0.030303030303030304
$ws.Range("M32").Value = 1
$ws.Range("N32").Value = 0.090909090909090912
$ws.Range("O32").Value = 3
$ws.Range("N37").Value = 0.29166666666666669
$ws.Range("O37").Value = 7
$ws.Range("D43").Value = 0.10526315789473684
$ws.Range("E43").Value = 10
$ws.Range("L43").Value = 0.10526315789473684
$ws.Range("M43").Value = 10
$ws.Range("D44").Value = 0.12
$ws.Range("E44").Value = 3
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0.12
$ws.Range("M44").Value = 3
$ws.Range("N44").Value = 0.2
$ws.Range("O44").Value = 5
$ws.Range("D45").Value = 0.083333333333333329
$ws.Range("E45").Value = 2
$ws.Range("J45").Value = 0.16666666666666666
$ws.Range("K45").Value = 4
$ws.Range("H46").Value = 0.31818181818181818
$ws.Range("I46").Value = 7
$ws.Range("J46").Value = 0.090909090909090912
$ws.Range("K46").Value = 2
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("F48").Value = 0.16666666666666666
$ws.Range("G48").Value = 5
$ws.Range("N49").Value = 0.1388888888888889
$ws.Range("O49").Value = 5
$ws.Range("D51").Value = 0.068181818181818177
$ws.Range("E51").Value = 3
$ws.Range("H51").Value = 0.20454545454545456
$ws.Range("I51").Value = 9
$ws.Range("F52").Value = 0.057142857142857141
$ws.Range("G52").Value = 2
$ws.Range("H54").Value = 0.078947368421052627
$ws.Range("I54").Value = 3
$ws.Range("L54").Value = 0.052631578947368418
$ws.Range("M54").Value = 2
